$wb = $excel.ActiveWorkbook

# --- Sheet: out_vars (sheet1) --- add row 26
$ws = $wb.Worksheets.Item("out_vars")
$ws.Range("A25:J25").Copy($ws.Range("A26:J26"))
$ws.Range("A26").Value = 44007
$ws.Range("B26").Value = 202951
$ws.Range("C26").Value = 262117
$ws.Range("D26").Value = 63583
$ws.Range("E26").Value = 25060
$ws.Range("F26").Value = 31.457839577040765
$ws.Range("G26").Value = 63844
$ws.Range("H26").Value = 5483
$ws.Range("I26").Value = 5681
$ws.Range("J26").Value = 528651

Write-Output "done sheet1"

# --- Sheet: dates_dx (sheet2) --- fill row 26
$ws = $wb.Worksheets.Item("dates_dx")
$ws.Range("A25").Copy($ws.Range("A26"))
$ws.Range("A26").Value = 44007
$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 1
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 4

$ws.Range("L26").Select()

Write-Output "done sheet2"

# --- Sheet: dates_sx (sheet3) --- fill row 26
$ws = $wb.Worksheets.Item("dates_sx")
$ws.Range("A25").Copy($ws.Range("A26"))
$ws.Range("A26").Value = 44007
$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 1
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 1
$ws.Range("K26").Value = 1
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 0

$ws.Range("D34").Select()

Write-Output "done sheet3"

# --- Sheet: dates_deaths (sheet4) --- fill row 26
$ws = $wb.Worksheets.Item("dates_deaths")
$ws.Range("A25").Copy($ws.Range("A26"))
$ws.Range("A26").Value = 44007
$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 1
$ws.Range("I26").Value = 1
$ws.Range("J26").Value = 2

$ws.Range("K26").Select()

Write-Output "done sheet4"

# --- Sheet: control_obs (sheet5) --- fill column Z (new date 6/25/20)
$ws = $wb.Worksheets.Item("control_obs")
$ws.Range("Y1").Copy($ws.Range("Z1"))
$ws.Range("Z1").Value = 44007
$ws.Range("Z2").Value = 3817
$ws.Range("Z3").Value = 3632
$ws.Range("Z4").Value = 3632
$ws.Range("Z5").Value = 3632
$ws.Range("Z6").Value = 3632
$ws.Range("Z7").Value = 2823
$ws.Range("Z8").Value = 5454
$ws.Range("Z10").Value = 165
$ws.Range("Z11").Value = 165
$ws.Range("Z12").Value = 165
$ws.Range("Z13").Value = 165
$ws.Range("Z14").Value = 165
$ws.Range("Z15").Value = 100
$ws.Range("Z16").Value = 177
$ws.Range("Z18").Value = 884
$ws.Range("Z20").Formula = "=SUM(Z2:Z18)"

$ws.Range("AA12").Select()

Write-Output "done sheet5"

# --- Sheet: out_vars (sheet1) --- becomes the active tab, with its own selection
$ws = $wb.Worksheets.Item("out_vars")
$ws.Activate()
$ws.Range("B23").Select()

Write-Output "done activating out_vars"
